$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.425.34'
Set-TextValue 'E2' '  +2.38%  '
Set-TextValue 'D3' '1.860.26'
Set-TextValue 'E3' '  +2.58%  '
Set-TextValue 'D4' '0.9972'
Set-TextValue 'E4' '  -0.54%  '
Set-TextValue 'D5' '281.34'
Set-TextValue 'E5' '  +0.72%  '
Set-TextValue 'D6' '0.9977'
Set-TextValue 'E6' '  -0.46%  '
Set-TextValue 'D7' '0.5127'
Set-TextValue 'E7' '  +2.29%  '
Set-TextValue 'D8' '0.3517'
Set-TextValue 'E8' '  -0.02%  '
Set-TextValue 'D9' '45.11'
Set-TextValue 'E9' '  +1.34%  '
Set-TextValue 'D10' '0.06849'
Set-TextValue 'E10' '  +2.88%  '
Set-TextValue 'D11' '20.04'
Set-TextValue 'E11' '  -0.04%  '
Set-TextValue 'D12' '0.8129'
Set-TextValue 'E12' '  -4.43%  '
Set-TextValue 'D13' '0.07759'
Set-TextValue 'E13' '  -0.96%  '
Set-TextValue 'D14' '1.858.17'
Set-TextValue 'E14' '  +2.55%  '
Set-TextValue 'D15' '89.03'
Set-TextValue 'E15' '  +1.82%  '
Set-TextValue 'D16' '5.105'
Set-TextValue 'E16' '  +1.35%  '
Set-TextValue 'D17' '0.9971'
Set-TextValue 'E17' '  -0.50%  '
Set-TextValue 'D18' '14.27'
Set-TextValue 'E18' '  +2.13%  '
Set-TextValue 'D19' '0.000008109'
Set-TextValue 'E19' '  +1.23%  '
Set-TextValue 'D20' '0.9972'
Set-TextValue 'E20' '  -0.54%  '
Set-TextValue 'D21' '26.433.76'
Set-TextValue 'E21' '  +2.22%  '
Set-TextValue 'D22' '4.790'
Set-TextValue 'E22' '  +0.78%  '
Set-TextValue 'D23' '10.10'
Set-TextValue 'E23' '  +1.18%  '
Set-TextValue 'D24' '6.217'
Set-TextValue 'E24' '  +2.08%  '
Set-TextValue 'D25' '2.366'
Set-TextValue 'E25' '  +10.77%  '
Set-TextValue 'D26' '144.36'
Set-TextValue 'E26' '  +1.85%  '
Set-TextValue 'B27' 'Toncoin'
Set-TextValue 'C27' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D27' '1.659'
Set-TextValue 'E27' '  -0.88%  '
Set-TextValue 'B28' 'EthereumClassic'
Set-TextValue 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D28' '17.28'
Set-TextValue 'E28' '  +2.81%  '
Set-TextValue 'D29' '110.45'
Set-TextValue 'E29' '  +1.41%  '
Set-TextValue 'D30' '4.375'
Set-TextValue 'E30' '  +1.96%  '
Set-TextValue 'D31' '4.318'
Set-TextValue 'E31' '  +2.16%  '
Set-TextValue 'D32' '0.08781'
Set-TextValue 'E32' '  -0.19%  '
Set-TextValue 'D33' '0.04900'
Set-TextValue 'E33' '  +1.91%  '
Set-TextValue 'D35' '0.7419'
Set-TextValue 'E35' '  +0.36%  '
Set-TextValue 'D36' '2.863'
Set-TextValue 'E36' '  +0.13%  '
Set-TextValue 'D37' '3.252'
Set-TextValue 'E37' '  +6.03%  '
Set-TextValue 'E38' '  -4.46%  '
Set-TextValue 'D39' '0.01861'
Set-TextValue 'E39' '  +0.38%  '
Set-TextValue 'D40' '0.5218'
Set-TextValue 'E40' '  -2.41%  '
Set-TextValue 'D41' '0.9628'
Set-TextValue 'E41' '  -1.12%  '
Set-TextValue 'D42' '116.28'
Set-TextValue 'E42' '  +3.32%  '
Set-TextValue 'D43' '6.274'
Set-TextValue 'E43' '  +1.13%  '
Set-TextValue 'D44' '8.038'
Set-TextValue 'E44' '  -1.92%  '
Set-TextValue 'D45' '0.9969'
Set-TextValue 'E45' '  -0.57%  '
Set-TextValue 'D46' '0.4552'
Set-TextValue 'E46' '  -3.41%  '
Set-TextValue 'D47' '0.1363'
Set-TextValue 'E47' '  -1.24%  '
Set-TextValue 'D48' '9.370'
Set-TextValue 'E48' '  +1.21%  '
Set-TextValue 'D49' '36.38'
Set-TextValue 'E49' '  +1.45%  '
Set-TextValue 'E50' '  +1.15%  '
Set-TextValue 'D51' '0.05917'
Set-TextValue 'E51' '  +0.33%  '
